# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.624.44"
$ws.Range("E2").Value = "  +2.42%  "
$ws.Range("D3").Value = "1.789.54"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.83"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.560"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "33.00"
$ws.Range("E8").Value = "  +7.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.282"
$ws.Range("E9").Value = "  +1.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0682"
$ws.Range("E10").Value = "  +3.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0938"
$ws.Range("E11").Value = "  +1.77%  "
$ws.Range("D12").Value = "2.044.59"
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.14"
$ws.Range("E13").Value = "  +11.31%  "
$ws.Range("D14").Value = "1.784.92"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("D16").Value = "34.558.90"
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.30"
$ws.Range("E17").Value = "  +2.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.64"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "254.01"
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("D20").Value = "0.0₃0776"
$ws.Range("E20").Value = "  +5.30%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.45"
$ws.Range("E22").Value = "  +1.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.25"
$ws.Range("E23").Value = "  +1.58%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.06"
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.37"
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.11"
$ws.Range("E27").Value = "  +2.45%  "
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.76"
$ws.Range("E30").Value = "  -0.94%  "
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.86"
$ws.Range("E34").Value = "  +2.36%  "
$ws.Range("D35").Value = "1.445.78"
$ws.Range("E35").Value = "  -2.53%  "
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0190"
$ws.Range("E37").Value = "  +2.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.630"
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "83.39"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("E40").Value = "  +4.83%  "
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.903"
$ws.Range("E42").Value = "  +1.88%  "
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0504"
$ws.Range("E44").Value = "  -2.05%  "
$ws.Range("E45").Value = "  +2.50%  "
$ws.Range("E46").Value = "  -2.43%  "
$ws.Range("D47").Value = "1.941.67"
$ws.Range("E47").Value = "  +0.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.92"
$ws.Range("E48").Value = "  +7.58%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.03"
$ws.Range("E49").Value = "  +2.47%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.39"
$ws.Range("E51").Value = "  -2.49%  "
